$wb = $excel.ActiveWorkbook

# --- Work on the "TimePeriods" sheet: add a new "Def7" column (H) ------------
$ws = $wb.Worksheets.Item("TimePeriods")

# Make this the active/selected sheet (mirrors the author moving to this tab
# and leaving the cursor near the bottom of the newly extended column).
$ws.Activate()

# Header cell H27 - copy the look of the neighbouring "Def6" header (G27) and
# give it the new label.
$ws.Range("G27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H27").Value = "Def7"

# Column H gets formatted like column G for the whole block of data rows
# (29 down to 79) before the numbers are typed in - this matches the way the
# existing D/G columns were built further up the sheet.
$ws.Range("G29").Copy() | Out-Null
$ws.Range("H29:H79").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new "Def7" time period values (2001 - 2050).
$values = @()
$year = 2001
for ($r = 29; $r -le 78; $r++) {
    $ws.Cells.Item($r, 8).Value = $year
    $year = $year + 1
}

# Widen the used range's view a bit and park the selection the way the author
# left it.
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q65").Select()

# --- "reporting options" sheet is no longer the tab shown on open -----------
# (handled automatically by activating "TimePeriods" above, which moves
# tabSelected there and clears it everywhere else)
